# Apply the odds updates from the 2024-10-17 FlashScore refresh.
# Each entry maps a cell reference on Sheet1 to its new numeric value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "M2" = 1.07
    "N2" = 9
    "O2" = 1.36
    "P2" = 3.2
    "Q2" = 2.08
    "R2" = 1.73
    "G3" = 1.9
    "H3" = 3.4
    "I3" = 4.33
    "J3" = 2.63
    "L3" = 5
    "X3" = 7.5
    "Z3" = 15
    "AE3" = 21
    "AF3" = 81
    "AH3" = 9.5
    "AI3" = 21
    "AK3" = 51
    "AZ3" = 101
    "BA3" = 151
    "H5" = 3.5
    "J5" = 3.5
    "M5" = 1.06
    "N5" = 10
    "X5" = 13
    "AC5" = 10
    "AG5" = 251
    "AH5" = 8
    "AU5" = 8
    "G6" = 1.38
    "H6" = 4.25
    "I6" = 6.6
    "J6" = 1.85
    "K6" = 2.37
    "L6" = 6.1
    "Q6" = 1.6
    "R6" = 2.05
    "U6" = 1.87
    "V6" = 1.89
    "X6" = 5.9
    "Z6" = 7.8
    "AA6" = 9.25
    "AD6" = 7.5
    "AE6" = 15
    "AF6" = 60
    "AH6" = 16
    "AI6" = 35
    "AJ6" = 17
    "AK6" = 110
    "AL6" = 55
    "AN6" = 3.25
    "AO6" = 6.2
    "AP6" = 15
    "AQ6" = 17
    "AR6" = 40
    "AT6" = 3.15
    "AW6" = 8.25
    "AX6" = 37
    "AY6" = 37
    "G7" = 2.62
    "I7" = 2.62
    "J7" = 3.15
    "K7" = 2.05
    "L7" = 3.2
    "Q7" = 1.72
    "R7" = 1.9
    "V7" = 2.2
    "W7" = 10.5
    "AA7" = 20
    "AB7" = 23
    "AD7" = 6.1
    "AH7" = 9.5
    "AP7" = 18
    "AQ7" = 55
    "AR7" = 75
    "AS7" = 175
    "AT7" = 2.77
    "AW7" = 4.7
    "AX7" = 14.5
    "AY7" = 19.5
    "AZ7" = 65
    "BA7" = 90
    "U8" = 1.69
    "H9" = 3.9
    "I9" = 6
    "K9" = 2.38
    "Q9" = 1.7
    "R9" = 2.1
    "S9" = 1.33
    "T9" = 3.25
    "U9" = 1.8
    "V9" = 1.8
    "W9" = 7.5
    "X9" = 7.5
    "AC9" = 12
    "AD9" = 7.5
    "AE9" = 17
    "AG9" = 251
    "AK9" = 67
    "AP9" = 19
    "AS9" = 126
    "AT9" = 3.25
    "AU9" = 8.5
    "BA9" = 126
    "M10" = 1.08
    "N10" = 8
    "Q10" = 2.3
    "R10" = 1.6
    "Q13" = 1.9
    "R13" = 1.95
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
